# Fruta / hortaliza, semanal
# Insert a new weekly record as row 5, pushing existing rows 5-13 down to 6-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 5 (shifts rows 5..13 -> 6..14)
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly data point
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44447
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107002
$ws.Range("J5").Value = "Chirimoya"
$ws.Range("K5").Value = "Cultivar IV Región"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 48
$ws.Range("N5").Value = 30000
$ws.Range("O5").Value = 30000
$ws.Range("P5").Value = 30000
$ws.Range("Q5").Value = "`$/bandeja 10 kilos"
$ws.Range("R5").Value = "Provincia del Elquí"
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 10

# Make sure the date cell keeps the workbook's date style (same as the other D-column cells)
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat

Write-Host "done"
